$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the paragraph "金牛犢, 金壇, 禁食" (the italic run directly
#    under the "jin" heading paragraph). Find it by its text and
#    delete the whole paragraph (including its paragraph mark).
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("金牛犢, 金壇, 禁食", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Paragraphs(1).Range.Delete()
}

# ------------------------------------------------------------------
# 2. Remove the paragraph "This PDF version is provided under the
#    same license."
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("This PDF version is provided under the same license.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Paragraphs(1).Range.Delete()
}

# ------------------------------------------------------------------
# 3. Rewrite the resource/license paragraph that starts with
#    "關鍵詞 (Biblica)" (bold) ... "CC BY-SA 4.0 license." (with two
#    hyperlinks inside). First drop the hyperlinks (keeping their
#    visible text as plain runs) so the whole span becomes ordinary
#    text, then replace the text and restore the bold formatting on
#    just the new title.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("關鍵詞 (Biblica) (Chinese (Traditional)) is based on", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $targetPara = $rng.Paragraphs(1)

    # Un-hyperlink every hyperlink that lives inside this paragraph
    # (there are exactly two: "Biblica, Inc." and "CC BY-SA 4.0
    # license"), from last to first so indices stay valid.
    for ($i = $d.Hyperlinks.Count; $i -ge 1; $i--) {
        $h = $d.Hyperlinks($i)
        if ($h.Range.Start -ge $targetPara.Range.Start -and $h.Range.End -le $targetPara.Range.End) {
            $h.Delete()
        }
    }

    $targetPara.Range.Find.Execute("關鍵詞 (Biblica) (Chinese (Traditional)) is based on: Biblica Bible Dictionary, Biblica, Inc., 2023, which is licensed under a CC BY-SA 4.0 license.", $true, $false, $false, $false, $false, $true, 1, $false, "Biblica Study Notes (Key Terms) © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. Biblica Study Notes has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual.", 2)

    # Restore non-bold formatting to everything after the new title.
    $titleRng = $d.Content
    $titleRng.Find.Execute("Biblica Study Notes (Key Terms)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $restRng = $d.Range($titleRng.End, $targetPara.Range.End)
    $restRng.Bold = 0
}

# ------------------------------------------------------------------
# 4. Remove the standalone "License Information" heading paragraph.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("License Information", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Paragraphs(1).Range.Delete()
}
